# Weekly data refresh: a new price observation is inserted as row 524
# (pushing the existing rows 524:561 down to 525:562) for
# "Femacal de La Calera - Poroto verde".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 524, shifting rows 524:561 down to 525:562.
$ws.Rows.Item(524).Insert()

# Populate the newly inserted row 524 with the new observation.
$ws.Range("A524").Value = 3
$ws.Range("B524").Value = "Femacal de La Calera"
$ws.Range("C524").Value = "Coquimbo"
$ws.Range("D524").Value = 45021
$ws.Range("E524").Value = 5
$ws.Range("F524").Value = 100112031
$ws.Range("G524").Value = "Poroto verde"
$ws.Range("H524").Value = "Magnum"
$ws.Range("I524").Value = "Primera"
$ws.Range("J524").Value = 73
$ws.Range("K524").Value = 18000
$ws.Range("L524").Value = 19000
$ws.Range("M524").Value = 18521
$ws.Range("N524").Value = "`$/saco 25 kilos"
$ws.Range("O524").Value = "Provincia de Quillota"
$ws.Range("P524").Value = 741
$ws.Range("Q524").Value = 25
$ws.Range("R524").Value = "Hortaliza"
